$d = $word.ActiveDocument

# Append a new bulleted paragraph (same list/number/font formatting as the
# rest of the document) with a red "commit code" note, ending with the
# Wingdings smiley symbol used elsewhere in the document.
$endRange = $d.Range($d.Content.End, $d.Content.End)

$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr>' +
  '<w:pStyle w:val="ListParagraph"/>' +
  '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
  '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="FF0000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
  '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="FF0000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
  '<w:t xml:space="preserve">Khi commit code, cần kiểm tra config , nếu bị thì phải lấy cái mới về , thêm phần mình vào , build chạy ok mới commit lên , vì đưa lên bị lỗi người khác rất khó tìm lỗi sữa lắm </w:t>' +
  '</w:r>' +
  '<w:r>' +
  '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="FF0000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' +
  '<w:sym w:font="Wingdings" w:char="F04A"/>' +
  '</w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$endRange.InsertXML($newParaXml)
